$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values would otherwise be auto-converted to numbers by Excel;
# temporarily force Text format so the value is stored as a string, matching the source data.
$protectedCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D15",
    "D17",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D47",
    "D48",
    "D50",
    "D51",
)
foreach ($addr in $protectedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.573.73"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.513.06"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "317.98"
$ws.Range("E5").Value = "  +4.50%  "
$ws.Range("D6").Value = "95.06"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "35.84"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "7.53"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("E13").Value = "  -2.27%  "
$ws.Range("D14").Value = "2.894.41"
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "15.49"
$ws.Range("E15").Value = "  +4.58%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.509.68"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  -2.88%  "
$ws.Range("D18").Value = "42.573.76"
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "0.0₃0967"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "6.53"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "71.23"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "250.99"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "2.99"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("D26").Value = "26.68"
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +12.49%  "
$ws.Range("D29").Value = "38.99"
$ws.Range("D30").Value = "10.04"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "155.85"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").Value = "3.33"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "2.07"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("D36").Value = "0.0783"
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("E37").Value = "  -5.26%  "
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "23.81"
$ws.Range("E40").Value = "  -7.56%  "
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "3.84"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "0.0300"
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").Value = "2.044.46"
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").Value = "84.26"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").Value = "8.80"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").Value = "2.752.69"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").Value = "72.53"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "0.189"
$ws.Range("E51").Value = "  -1.26%  "

# Restore the default (Normal) style on the protected cells now that the text value is set,
# so no stray number-format styling is left behind on them.
foreach ($addr in $protectedCells) {
    $ws.Range($addr).Style = "Normal"
}
